# EPEX Spot prices workbook update:
# Add a new "18-jun" price column (E) to the "Prix Spot" sheet -
# a header cell styled like the existing day headers, and a "-"
# placeholder (no data published yet) for every hourly row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Prix Spot")

# Copy the formatting of the previous day's header cell (D1) onto the
# new header cell (E1) so it keeps the bold/centered/bordered header
# style, then set its text.
$ws.Range("D1").Copy()
$ws.Range("E1").PasteSpecial(-4122)
$ws.Range("E1").Value = "18-jun"

# Fill the new column's data rows (2-25) with the "-" placeholder used
# for not-yet-available prices.
for ($row = 2; $row -le 25; $row++) {
    $ws.Cells.Item($row, 5).Value = "-"
}
